$wb = $excel.ActiveWorkbook

# Define the text replacements that need to happen inside the relevant cells.
# Order matters: longer / more specific patterns first so we don't double-replace.
function Update-ScenarioText($text) {
    if ($null -eq $text) { return $text }
    $new = $text
    $new = $new.Replace("Scenario 4 - 100% Cost Share", "Scenario 4 - 100% Capex`nShare")
    $new = $new.Replace("Scenario 4 - 30% Cost Share", "Scenario 4 - 30% Capex`nShare")
    $new = $new.Replace("Scenario 4 - 50% Cost Share", "Scenario 4 - 50% Capex`nShare")
    $new = $new.Replace("Scenario 4 - Low Rate", "Scenario 4 - Low Elec.`nRate")
    return $new
}

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    $cols = $used.Columns.Count
    for ($r = 1; $r -le $rows; $r++) {
        for ($c = 1; $c -le $cols; $c++) {
            $cell = $used.Cells.Item($r, $c)
            $val = $cell.Value2
            if ($val -is [string]) {
                $updated = Update-ScenarioText $val
                if ($updated -ne $val) {
                    $cell.Value = $updated
                }
            }
        }
    }
}
